# add c docs (in code comments)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("Feuil2")

# --- Feuil1: append new data rows (12-19) ---
$data1 = @(
    @(-6, -2),
    @(-7, 5),
    @(2, -8),
    @(1, -2),
    @(9, -3),
    @(5, -4),
    @(-3, 7),
    @(-2, 10)
)

$r = 12
foreach ($row in $data1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# --- Feuil2: update existing values in rows 2-3 ---
$ws2.Range("A2").Value = 5
$ws2.Range("B2").Value = 80
$ws2.Range("C2").Value = 200

$ws2.Range("A3").Value = 4
$ws2.Range("B3").Value = 50
$ws2.Range("C3").Value = 200

# --- Selection state: Feuil2 goes back to its original B2 selection ---
$ws2.Select()
$ws2.Range("B2").Select()

# --- Feuil1 becomes the active / tab-selected sheet, selection A3:B3 ---
$ws1.Select()
$ws1.Range("A3:B3").Select()
